$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.860.65'
$ws.Range('E2').Value = '  +0.90%  '
$ws.Range('D3').Value = '3.542.83'
$ws.Range('E3').Value = '  -0.05%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '615.82'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.26%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '152.89'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.41%  '
$ws.Range('D7').Value = '3.540.90'
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.481'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.16%  '
$ws.Range('E10').Value = '  -1.17%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.09'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.80%  '
$ws.Range('E12').Value = '  -1.22%  '
$ws.Range('E13').Value = '  -1.09%  '
$ws.Range('D14').Value = '4.148.19'
$ws.Range('E14').Value = '  +0.17%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '32.17'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.23%  '
$ws.Range('D16').Value = '3.555.93'
$ws.Range('E16').Value = '  +0.81%  '
$ws.Range('D17').Value = '67.668.71'
$ws.Range('E17').Value = '  +0.68%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.117'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.63%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.41'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.01%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.36'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.93%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '447.86'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.38%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.67'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.46%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.625'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '77.56'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.35%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000133'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +6.18%  '
$ws.Range('D26').Value = '3.686.56'
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('B27').Value = 'Dai'
$ws.Range('C27').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.05%  '
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.29'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.54%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.64'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.81%  '
$ws.Range('E30').Value = '  -1.52%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.61'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.71%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.169'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +6.41%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '25.98'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.20%  '
$ws.Range('E35').Value = '  -0.09%  '
$ws.Range('D36').Value = '3.532.60'
$ws.Range('E36').Value = '  -0.15%  '
$ws.Range('E37').Value = '  -2.41%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.05'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.48%  '
$ws.Range('E39').Value = '  -0.04%  '
$ws.Range('E40').Value = '  +0.03%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '176.84'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.18%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0897'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.61%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.19'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.77%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.43'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.85%  '
$ws.Range('E45').Value = '  -0.74%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '28.58'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.03%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '45.39'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.97%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.68'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.97%  '
$ws.Range('E49').Value = '  +4.21%  '
$ws.Range('E50').Value = '  -0.74%  '
$ws.Range('E51').Value = '  -3.27%  '
